$d = $word.ActiveDocument

# The review body (paragraphs 2-10) is collapsed into a single short
# paragraph reading "Inoltre ". We keep paragraph 1 (title) and the two
# trailing empty paragraphs untouched.

# Insert a brand-new (plain, unformatted) paragraph right before the big
# block of review text, so the replacement paragraph doesn't inherit the
# "La struttura dati..." paragraph's underline-on-mark formatting.
$pFirstReview = $d.Paragraphs.Item(2)
$pFirstReview.Range.InsertParagraphBefore()

# After the insert, the old paragraph 2 ("Game component ...") is now
# paragraph 3, and the old paragraph 10 ("La struttura dati ...") is now
# paragraph 11. Delete that whole span (content + paragraph marks) in one
# shot, collapsing it away entirely.
$pBlockStart = $d.Paragraphs.Item(3)
$blockEnd = $d.Paragraphs.Item(11).Range.End
$blockRange = $d.Range($pBlockStart.Range.Start, $blockEnd)
$blockRange.Delete()

# Fill the freshly inserted (now second) paragraph with the replacement text.
$pNew = $d.Paragraphs.Item(2)
$pNew.Range.Text = "Inoltre "
